$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$text)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = "55.121.94"
$ws.Range("D3").Value = "2.343.07"
$ws.Range("E3").Value = "  -5.43%  "
$ws.Range("E4").Value = "  -0.07%  "
Set-TextValue $ws.Range("D5") "475.14"
$ws.Range("E5").Value = "  -2.74%  "
Set-TextValue $ws.Range("D6") "145.57"
$ws.Range("E6").Value = "  -0.87%  "
Set-TextValue $ws.Range("D7") "0.632"
$ws.Range("E7").Value = "  +24.22%  "
Set-TextValue $ws.Range("D8") "0.998"
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "2.348.03"
$ws.Range("E9").Value = "  -5.57%  "
Set-TextValue $ws.Range("D10") "0.0961"
$ws.Range("E10").Value = "  -0.84%  "
Set-TextValue $ws.Range("D11") "5.44"
$ws.Range("E11").Value = "  -6.44%  "
Set-TextValue $ws.Range("D12") "0.324"
$ws.Range("E12").Value = "  -2.27%  "
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("D14").Value = "2.749.31"
$ws.Range("E14").Value = "  -5.48%  "
$ws.Range("D15").Value = "55.092.44"
$ws.Range("E15").Value = "  -1.97%  "
Set-TextValue $ws.Range("D16") "19.95"
$ws.Range("E16").Value = "  -5.46%  "
Set-TextValue $ws.Range("D17") "0.0000129"
$ws.Range("E17").Value = "  -4.89%  "
$ws.Range("D18").Value = "2.342.86"
$ws.Range("E18").Value = "  -5.70%  "
$ws.Range("E19").Value = "  +0.92%  "
Set-TextValue $ws.Range("D20") "314.77"
$ws.Range("E20").Value = "  -1.07%  "
Set-TextValue $ws.Range("D21") "9.60"
$ws.Range("E21").Value = "  -4.74%  "
Set-TextValue $ws.Range("D22") "0.998"
$ws.Range("E22").Value = "  -0.05%  "
Set-TextValue $ws.Range("D23") "5.62"
$ws.Range("E23").Value = "  -3.28%  "
Set-TextValue $ws.Range("D24") "56.76"
$ws.Range("E24").Value = "  -2.81%  "
$ws.Range("E25").Value = "  +0.00%  "
Set-TextValue $ws.Range("D26") "0.394"
$ws.Range("E26").Value = "  -4.49%  "
$ws.Range("E27").Value = "  -6.73%  "
$ws.Range("D28").Value = "2.442.72"
$ws.Range("E28").Value = "  -5.35%  "
Set-TextValue $ws.Range("D29") "7.01"
$ws.Range("E29").Value = "  -7.98%  "
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("D31").Value = "0.0₃0743"
$ws.Range("E31").Value = "  -6.04%  "
$ws.Range("E32").Value = "  -0.29%  "
Set-TextValue $ws.Range("D33") "143.90"
$ws.Range("E33").Value = "  -3.50%  "
$ws.Range("E34").Value = "  -2.39%  "
Set-TextValue $ws.Range("D35") "5.12"
$ws.Range("E35").Value = "  -1.68%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D36") "1.09"
$ws.Range("E36").Value = "  -4.76%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D37") "3.57"
$ws.Range("E37").Value = "  -4.60%  "
Set-TextValue $ws.Range("D38") "0.808"
$ws.Range("E38").Value = "  -6.21%  "
Set-TextValue $ws.Range("D39") "0.102"
$ws.Range("E39").Value = "  +10.01%  "
Set-TextValue $ws.Range("D40") "33.68"
$ws.Range("E40").Value = "  -1.53%  "
Set-TextValue $ws.Range("D41") "0.998"
$ws.Range("E41").Value = "  +0.23%  "
Set-TextValue $ws.Range("D42") "3.41"
$ws.Range("E42").Value = "  -3.06%  "
$ws.Range("E43").Value = "  -0.64%  "
Set-TextValue $ws.Range("D44") "0.575"
$ws.Range("E44").Value = "  -4.96%  "
Set-TextValue $ws.Range("D45") "0.0517"
$ws.Range("E45").Value = "  -7.08%  "
Set-TextValue $ws.Range("D46") "10.17"
$ws.Range("E46").Value = "  -0.31%  "
Set-TextValue $ws.Range("D47") "250.02"
$ws.Range("E47").Value = "  -2.95%  "
$ws.Range("E48").Value = "  -3.86%  "
Set-TextValue $ws.Range("D49") "4.31"
$ws.Range("E49").Value = "  -9.42%  "
Set-TextValue $ws.Range("D50") "16.64"
$ws.Range("E50").Value = "  -5.59%  "
$ws.Range("D51").Value = "1.780.17"
$ws.Range("E51").Value = "  -5.01%  "
